# Update balance-sheet database: shift yearly columns and refresh headers/dates
$ws = $excel.ActiveSheet

# --- Header: financial period labels (row 8) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Header: publish dates (row 9) ---
$ws.Range("D9").Value = "1399-02-09 (10)"
$ws.Range("E9").Value = "1400-02-01 (11)"
$ws.Range("F9").Value = "1401-02-10 (12)"
$ws.Range("G9").Value = "1402-02-09 (10)"
$ws.Range("H9").Value = "1402-02-09 (2)"

# --- Balance sheet data: shift columns left (drop oldest period, add newest) ---
# Row 12
$ws.Range("D12").Value = 19416
$ws.Range("E12").Value = 53910
$ws.Range("F12").Value = 112105
$ws.Range("G12").Value = 83775
$ws.Range("H12").Value = 380934
# Row 13
$ws.Range("D13").Value = 179883
$ws.Range("E13").Value = 447490
$ws.Range("F13").Value = 323319
$ws.Range("G13").Value = 78583
$ws.Range("H13").Value = 24818
# Row 14
$ws.Range("D14").Value = 273863
$ws.Range("E14").Value = 242903
$ws.Range("F14").Value = 464557
$ws.Range("G14").Value = 801486
$ws.Range("H14").Value = 1242052
# Row 15
$ws.Range("D15").Value = 271850
$ws.Range("E15").Value = 415305
$ws.Range("F15").Value = 584064
$ws.Range("G15").Value = 662112
$ws.Range("H15").Value = 982961
# Row 16
$ws.Range("D16").Value = 28856
$ws.Range("E16").Value = 49477
$ws.Range("F16").Value = 29960
$ws.Range("G16").Value = 67273
$ws.Range("H16").Value = 31850
# Row 18
$ws.Range("D18").Value = 773868
$ws.Range("E18").Value = 1209085
$ws.Range("F18").Value = 1514005
$ws.Range("G18").Value = 1693229
$ws.Range("H18").Value = 2662615
# Row 20
$ws.Range("D20").Value = 64269
$ws.Range("E20").Value = 64270
$ws.Range("F20").Value = 64269
$ws.Range("G20").Value = 80623
$ws.Range("H20").Value = 80623
# Row 22
$ws.Range("D22").Value = 57333
$ws.Range("E22").Value = 234618
$ws.Range("F22").Value = 933812
$ws.Range("G22").Value = 1885628
$ws.Range("H22").Value = 2134428
# Row 23
$ws.Range("D23").Value = 4356
$ws.Range("E23").Value = 5053
$ws.Range("F23").Value = 5121
$ws.Range("G23").Value = 6564
$ws.Range("H23").Value = 6578
# Row 25
$ws.Range("D25").Value = 42380
$ws.Range("E25").Value = 31229
$ws.Range("F25").Value = 27594
$ws.Range("G25").Value = 10003
$ws.Range("H25").Value = 50735
# Row 26
$ws.Range("D26").Value = 168338
$ws.Range("E26").Value = 335170
$ws.Range("F26").Value = 1030796
$ws.Range("G26").Value = 1982818
$ws.Range("H26").Value = 2272364
# Row 27
$ws.Range("D27").Value = 942206
$ws.Range("E27").Value = 1544255
$ws.Range("F27").Value = 2544801
$ws.Range("G27").Value = 3676047
$ws.Range("H27").Value = 4934979
# Row 29
$ws.Range("D29").Value = 252318
$ws.Range("E29").Value = 234660
$ws.Range("F29").Value = 311786
$ws.Range("G29").Value = 586802
$ws.Range("H29").Value = 769828
# Row 31
$ws.Range("D31").Value = 3150
$ws.Range("E31").Value = 37565
$ws.Range("F31").Value = 23512
$ws.Range("G31").Value = 33081
$ws.Range("H31").Value = 99511
# Row 32
$ws.Range("D32").Value = 23984
$ws.Range("E32").Value = 81126
$ws.Range("F32").Value = 145950
$ws.Range("G32").Value = 216161
$ws.Range("H32").Value = 250566
# Row 33
$ws.Range("D33").Value = 28807
$ws.Range("E33").Value = 26882
$ws.Range("F33").Value = 58409
$ws.Range("G33").Value = 74740
$ws.Range("H33").Value = 88121
# Row 34
$ws.Range("D34").Value = 55502
$ws.Range("E34").Value = 14290
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 41667
# Row 37
$ws.Range("D37").Value = 363761
$ws.Range("E37").Value = 394523
$ws.Range("F37").Value = 539657
$ws.Range("G37").Value = 910784
$ws.Range("H37").Value = 1249693
# Row 38
$ws.Range("D38").Value = 32483
$ws.Range("E38").Value = 40455
$ws.Range("F38").Value = 46241
$ws.Range("G38").Value = 57869
$ws.Range("H38").Value = 78477
# Row 39
$ws.Range("D39").Value = "-"
# Row 40
$ws.Range("D40").Value = 18296
$ws.Range("E40").Value = 4871
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0
# Row 41
$ws.Range("D41").Value = 78344
$ws.Range("E41").Value = 114337
$ws.Range("F41").Value = 127587
$ws.Range("G41").Value = 177906
$ws.Range("H41").Value = 257619
# Row 42
$ws.Range("D42").Value = 129123
$ws.Range("E42").Value = 159663
$ws.Range("F42").Value = 173828
$ws.Range("G42").Value = 235775
$ws.Range("H42").Value = 336096
# Row 43
$ws.Range("D43").Value = 492884
$ws.Range("E43").Value = 554186
$ws.Range("F43").Value = 713485
$ws.Range("G43").Value = 1146559
$ws.Range("H43").Value = 1585789
# Row 45
$ws.Range("D45").Value = 450000
$ws.Range("G45").Value = 804000
$ws.Range("H45").Value = 2000000
# Row 47
$ws.Range("D47").Value = 0
$ws.Range("F47").Value = 354000
$ws.Range("G47").Value = 0
# Row 49
$ws.Range("D49").Value = 0
# Row 50
$ws.Range("E50").Value = 43037
$ws.Range("F50").Value = 45000
$ws.Range("G50").Value = 80400
$ws.Range("H50").Value = 132193
# Row 52
$ws.Range("D52").Value = "-"
# Row 54
$ws.Range("D54").Value = "-"
# Row 56
$ws.Range("D56").Value = -16678
$ws.Range("E56").Value = 497032
$ws.Range("F56").Value = 982316
$ws.Range("G56").Value = 1645088
$ws.Range("H56").Value = 1216997
# Row 57
$ws.Range("D57").Value = 449322
$ws.Range("E57").Value = 990069
$ws.Range("F57").Value = 1831316
$ws.Range("G57").Value = 2529488
$ws.Range("H57").Value = 3349190
# Row 58
$ws.Range("D58").Value = 942206
$ws.Range("E58").Value = 1544255
$ws.Range("F58").Value = 2544801
$ws.Range("G58").Value = 3676047
$ws.Range("H58").Value = 4934979
